$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.960.36"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.257.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.26%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.50%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.34%  "

# Row 9
$ws.Range("E9").Value = "  +3.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.415"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.49%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.825.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("E13").Value = "  +0.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.989.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.42%  "

# Row 16
$ws.Range("E16").Value = "  +2.50%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.256.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.91%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "381.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.75%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.514"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.34%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("E27").Value = "  +2.27%  "

# Row 28
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.30%  "

# Row 32
$ws.Range("E32").Value = "  +1.14%  "

# Row 33
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("E34").Value = "  +0.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.68%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.14%  "

# Row 37
$ws.Range("E37").Value = "  -0.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.834"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.98%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.46%  "

# Row 40
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.21%  "

# Row 42
$ws.Range("E42").Value = "  +1.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.82%  "

# Row 45
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0687"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.79%  "

# Row 46
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "345.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.641.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.58%  "

# Row 48
$ws.Range("E48").Value = "  +1.18%  "

# Row 49
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.64%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "

# Row 51
$ws.Range("E51").Value = "  +1.02%  "
